# Update "想去人数" (F column) figures on the 展览 (sheet index 1) and
# 全部类型 (sheet index 4) worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 306
$ws1.Range("F3").Value = 1230
$ws1.Range("F4").Value = 16955
$ws1.Range("F6").Value = 1658
$ws1.Range("F7").Value = 71
$ws1.Range("F8").Value = 9
$ws1.Range("F9").Value = 9
$ws1.Range("F13").Value = 11742
$ws1.Range("F15").Value = 8
$ws1.Range("F16").Value = 1429
$ws1.Range("F17").Value = 4664
$ws1.Range("F18").Value = 474
$ws1.Range("F19").Value = 15
$ws1.Range("F22").Value = 904
$ws1.Range("F25").Value = 30

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 306
$ws4.Range("F4").Value = 1230
$ws4.Range("F5").Value = 16955
$ws4.Range("F7").Value = 1658
$ws4.Range("F8").Value = 71
$ws4.Range("F9").Value = 9
$ws4.Range("F10").Value = 9
$ws4.Range("F16").Value = 11742
$ws4.Range("F18").Value = 8
$ws4.Range("F19").Value = 1429
$ws4.Range("F20").Value = 4664
$ws4.Range("F21").Value = 474
$ws4.Range("F22").Value = 15
$ws4.Range("F25").Value = 904
$ws4.Range("F28").Value = 30
